$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 480031.68
$ws.Range("M4").Value = 480031.68
$ws.Range("H5").Value = 480031.68
$ws.Range("N5").Value = 480031.68

$ws.Range("G6").Value = 4053505.25

$ws.Range("M8").Value = 12695902.02

$ws.Range("H10").Value = 4053505.25
$ws.Range("N10").Value = 4053505.25
$ws.Range("O10").Value = 0

$ws.Range("G11").Value = 3236811.88

$ws.Range("M13").Value = 3236811.88

$ws.Range("H15").Value = 3236811.88
$ws.Range("N15").Value = 3236811.88
$ws.Range("O15").Value = 0

$ws.Range("G16").Value = 387874.28

$ws.Range("M18").Value = 100559.88

$ws.Range("M20").Value = 287314.4

$ws.Range("H21").Value = 387874.28
$ws.Range("N21").Value = 387874.28
$ws.Range("O21").Value = 0

$ws.Range("G25").Value = 1192366.06

$ws.Range("M27").Value = 1192366.06

$ws.Range("H28").Value = 1192366.06
$ws.Range("N28").Value = 1192366.06

$ws.Range("G33").Value = 251338.26

$ws.Range("M35").Value = 251338.26

$ws.Range("H36").Value = 251338.26
$ws.Range("N36").Value = 251338.26

$ws.Range("G41").Value = 701458.4

$ws.Range("M43").Value = 701458.4

$ws.Range("H44").Value = 701458.4
$ws.Range("N44").Value = 701458.4
